# Refresh cryptos list: latest Price/Volume(1h) scrape.
# Row 46/47 also swap (Maker <-> InjectiveProtocol reordered by rank).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.419.06"
$ws.Range("E2").Value = "  +0.62%  "

$ws.Range("D3").Value = "3.472.04"

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.61"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +0.18%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.32"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +3.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.617"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +5.20%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").Value = "3.468.86"
$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.140"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +8.47%  "

$ws.Range("E11").Value = "  -1.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.429"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +0.81%  "

$ws.Range("D13").Value = "4.067.47"
$ws.Range("E13").Value = "  -0.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.97"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +3.96%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.134"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -0.48%  "

$ws.Range("D16").Value = "67.426.62"
$ws.Range("E16").Value = "  +0.59%  "

$ws.Range("E17").Value = "  +1.04%  "

$ws.Range("D18").Value = "3.468.61"
$ws.Range("E18").Value = "  +0.28%  "

$ws.Range("E19").Value = "  -0.59%  "

$ws.Range("E20").Value = "  -1.58%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "394.20"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +1.88%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.94"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +1.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.79"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +1.42%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +0.05%  "

$ws.Range("E25").Value = "  +1.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.81"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -0.94%  "

$ws.Range("E27").Value = "  +0.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.38"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +1.05%  "

$ws.Range("E29").Value = "  -1.19%  "

$ws.Range("E30").Value = "  +0.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.13"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +0.68%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.41"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -0.55%  "

$ws.Range("E33").Value = "  +1.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.54"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +0.44%  "

$ws.Range("E35").Value = "  +0.90%  "

$ws.Range("E36").Value = "  -0.15%  "

$ws.Range("E37").Value = "  -2.29%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.99"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -1.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.893"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +2.54%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.85"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +10.69%  "

$ws.Range("E41").Value = "  -2.66%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.78"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -2.33%  "

$ws.Range("E43").Value = "  +1.17%  "

$ws.Range("E44").Value = "  -0.21%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0718"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -0.46%  "

$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "26.27"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -3.86%  "

$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.740.33"
$ws.Range("E47").Value = "  -1.68%  "

$ws.Range("E48").Value = "  -1.37%  "

$ws.Range("E49").Value = "  +0.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "325.87"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -3.78%  "

$ws.Range("E51").Value = "  -2.35%  "

